$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename keyword cells to their longer forms
$ws.Range("A35").Value = "핫넘버"
$ws.Range("A36").Value = "콜드넘버"

# Reflect the selection change recorded for this sheet (user last selected F37)
$ws.Range("F37").Select()
